# MEC-3B schedule fix: shift the afternoon block down by one slot so the
# lunch break ("Almoço") starts at 12:20 (giving a full 6-hour span per
# turn) and add three new trailing time slots (16:40, 17:30, 18:20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 11:30 slot becomes a plain "-" row (used to hold "Almoço").
$ws.Range("A8").Value = "11:30"
$ws.Range("B8:F8").Value = "-"

# Row 9: new 12:20 slot becomes the "Almoço" row.
$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

# Row 10: shifts from 13:50 -> 13:00, stays "-".
$ws.Range("A10").Value = "13:00"
$ws.Range("B10:F10").Value = "-"

# Row 11: shifts from 14:40 -> 13:50, stays "-".
$ws.Range("A11").Value = "13:50"
$ws.Range("B11:F11").Value = "-"

# Row 12: shifts from 15:30 -> 14:40, and content becomes "-" (used to be
# "Intervalo").
$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"

# Row 13: shifts from 15:50 -> 15:30, content becomes "Intervalo".
$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

# Row 14: now holds what used to be row 13's time (15:50), content "-".
$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

# Row 15 (new): 16:40, "-".
$ws.Range("A15").Value = "16:40"
$ws.Range("B15:F15").Value = "-"

# Row 16 (new): 17:30, "-".
$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"

# Row 17 (new): 18:20, with the remaining cells left blank.
$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Value = ""
